# Update functions & data on Sheet1 (data/deaths_location.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Refresh the row 2 (Scotland) figures with the latest reported totals
$ws.Range("B2").Value = 44003
$ws.Range("C2").Value = 4119
$ws.Range("D2").Value = 1909
$ws.Range("E2").Value = 1917
$ws.Range("F2").Value = 286
$ws.Range("G2").Value = 7

# Move the active selection to J2, matching the saved cursor position
$ws.Activate()
$ws.Range("J2").Select()
